$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'51.190.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.45%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.949.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.69%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'378.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.63%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'104.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.21%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.542"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.18%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.06%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  +0.38%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'36.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.36%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +0.76%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  +0.75%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'3.421.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.66%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'18.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.13%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'7.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.24%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'2.951.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.88%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.952"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.05%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'51.252.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.67%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  +1.15%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = "'  +1.35%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.54%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  +1.64%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'68.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.94%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'260.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.50%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +4.10%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  -1.10%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = "'RenderToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'7.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +17.51%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "'Filecoin"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'7.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.49%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  +0.08%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("B30").Value = "'Hedera"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'0.112"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +8.93%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("B31").Value = "'EthereumClassic"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'25.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.00%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'9.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.82%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'34.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.36%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("B34").Value = "'Toncoin"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'2.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.99%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("B35").Value = "'OKB"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'51.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.45%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.0445"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.92%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  -0.20%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  -0.77%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'17.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.67%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  -4.03%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'1.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.65%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +2.17%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'123.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +5.07%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'22.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.73%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("B45").Value = "'TheGraph"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.280"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +16.65%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'WEMIXToken"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.49%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.92%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'2.034.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.42%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +0.85%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0348"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +11.22%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'5.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.16%  "
$ws.Range("E51").Style = "Normal"
